# Update crypto price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.709.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.55%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.426.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.03%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.73%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.27%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.418.32"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.90%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("E9").Value = "  -0.45%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.201"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.04%  "

# Row 11
$ws.Range("E11").Value = "  -0.96%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.10%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000283"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.24%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "694.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.75%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.971.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.83%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.09%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.719.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.48%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.425.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.06%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.122"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.08%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.67"
$ws.Range("D20").Style = "Normal"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.48%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.897"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.67%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.21%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.44%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.49%  "

# Row 26
$ws.Range("E26").Value = "  -0.48%  "

# Row 27
$ws.Range("E27").Value = "  -2.85%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.42%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.31%  "

# Row 30
$ws.Range("E30").Value = "  +0.14%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.89%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "572.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.15%  "

# Row 33
$ws.Range("E33").Value = "  +0.15%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.93%  "

# Row 35
$ws.Range("E35").Value = "  -2.60%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.16%  "

# Row 37
$ws.Range("E37").Value = "  +0.10%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.560.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.01%  "

# Row 39
$ws.Range("E39").Value = "  -0.71%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.17%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0735"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.71%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.19%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.09%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.52%  "

# Row 45
$ws.Range("E45").Value = "  -2.23%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0417"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.44%  "

# Row 47
$ws.Range("E47").Value = "  +3.79%  "

# Row 48
$ws.Range("E48").Value = "  -0.66%  "

# Row 49
$ws.Range("E49").Value = "  -1.18%  "

# Row 50
$ws.Range("E50").Value = "  -0.33%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.83%  "
